$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 0.09425133333333334
$ws.Range("H2").Value = 0.282754
$ws.Range("I2").Value = 0.02715992817009031
$ws.Range("J2").Value = 0.02715992817009031
$ws.Range("M2").Value = 2.325008666666667
$ws.Range("N2").Value = 6.975026
$ws.Range("O2").Value = 0.05445297772988467
$ws.Range("P2").Value = 0.05445297772988466
$ws.Range("Q2").Value = 0.2191351668448889
$ws.Range("R2").Value = 1.972216501604
$ws.Range("S2").Value = 0.001478938963791195
$ws.Range("T2").Value = 0.001478938963791194

$ws.Range("G3").Value = 0.09425133333333334
$ws.Range("H3").Value = 0.282754
$ws.Range("I3").Value = 0.02715992817009031
$ws.Range("J3").Value = 0.02715992817009031
$ws.Range("O3").Value = 0.4529132218878514
$ws.Range("P3").Value = 0.4529132218878514
$ws.Range("Q3").Value = 1.822659082795778
$ws.Range("R3").Value = 16.403931745162
$ws.Range("S3").Value = 0.01230109057375822
$ws.Range("T3").Value = 0.01230109057375822

$ws.Range("G4").Value = 0.09425133333333334
$ws.Range("H4").Value = 0.282754
$ws.Range("I4").Value = 0.02715992817009031
$ws.Range("J4").Value = 0.02715992817009031
$ws.Range("O4").Value = 0.492633800382264
$ws.Range("P4").Value = 0.492633800382264
$ws.Range("Q4").Value = 1.982506642257556
$ws.Range("R4").Value = 17.842559780318
$ws.Range("S4").Value = 0.0133798986325409
$ws.Range("T4").Value = 0.0133798986325409

$ws.Range("I5").Value = 0.357039508851706
$ws.Range("J5").Value = 0.357039508851706
$ws.Range("M5").Value = 2.325008666666667
$ws.Range("N5").Value = 6.975026
$ws.Range("O5").Value = 0.05445297772988467
$ws.Range("P5").Value = 0.05445297772988466
$ws.Range("Q5").Value = 2.880711313095333
$ws.Range("R5").Value = 25.926401817858
$ws.Range("S5").Value = 0.01944186442419091
$ws.Range("T5").Value = 0.0194418644241909

$ws.Range("I6").Value = 0.357039508851706
$ws.Range("J6").Value = 0.357039508851706
$ws.Range("O6").Value = 0.4529132218878514
$ws.Range("P6").Value = 0.4529132218878514
$ws.Range("S6").Value = 0.1617079142952822
$ws.Range("T6").Value = 0.1617079142952822

$ws.Range("I7").Value = 0.357039508851706
$ws.Range("J7").Value = 0.357039508851706
$ws.Range("O7").Value = 0.492633800382264
$ws.Range("P7").Value = 0.492633800382264
$ws.Range("S7").Value = 0.1758897301322329
$ws.Range("T7").Value = 0.1758897301322329

$ws.Range("I8").Value = 0.6158005629782037
$ws.Range("J8").Value = 0.6158005629782037
$ws.Range("M8").Value = 2.325008666666667
$ws.Range("N8").Value = 6.975026
$ws.Range("O8").Value = 0.05445297772988467
$ws.Range("P8").Value = 0.05445297772988466
$ws.Range("Q8").Value = 4.968479970429779
$ws.Range("R8").Value = 44.716319733868
$ws.Range("S8").Value = 0.03353217434190257
$ws.Range("T8").Value = 0.03353217434190257

$ws.Range("I9").Value = 0.6158005629782037
$ws.Range("J9").Value = 0.6158005629782037
$ws.Range("O9").Value = 0.4529132218878514
$ws.Range("P9").Value = 0.4529132218878514
$ws.Range("S9").Value = 0.278904217018811
$ws.Range("T9").Value = 0.278904217018811

$ws.Range("I10").Value = 0.6158005629782037
$ws.Range("J10").Value = 0.6158005629782037
$ws.Range("O10").Value = 0.492633800382264
$ws.Range("P10").Value = 0.492633800382264
$ws.Range("S10").Value = 0.3033641716174902
$ws.Range("T10").Value = 0.3033641716174902
